$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)  # ALC
$ws.Cells.Item(44, 8).Value = 15000  # H44: 0 -> 15000
$ws.Cells.Item(44, 10).Value = 15000  # J44: 0 -> 15000
$ws.Cells.Item(44, 12).Value = 15000  # L44: 0 -> 15000
$ws.Cells.Item(44, 14).Value = -15924  # N44: None -> -15924

$ws.Cells.Item(64, 8).Value = 6000.5  # H64: 4357.4287 -> 6000.5
$ws.Cells.Item(64, 9).Value = 4000.5  # I64: 3375 -> 4000.5
$ws.Cells.Item(64, 10).Value = 7000.5  # J64: 5667.3335 -> 7000.5
$ws.Cells.Item(64, 11).Value = 4000.5  # K64: 3375 -> 4000.5
$ws.Cells.Item(64, 12).Value = 7000.5  # L64: 5667.3335 -> 7000.5
$ws.Cells.Item(64, 13).Value = -3752.5  # M64: -3127 -> -3752.5
$ws.Cells.Item(64, 14).Value = -7496.5  # N64: -6163.3335 -> -7496.5

$ws.Cells.Item(67, 8).Value = 6000.5  # H67: 4357.4287 -> 6000.5
$ws.Cells.Item(67, 9).Value = 4000.5  # I67: 3375 -> 4000.5
$ws.Cells.Item(67, 10).Value = 7000.5  # J67: 5667.3335 -> 7000.5
$ws.Cells.Item(67, 11).Value = 4000.5  # K67: 3375 -> 4000.5
$ws.Cells.Item(67, 12).Value = 7000.5  # L67: 5667.3335 -> 7000.5
$ws.Cells.Item(67, 13).Value = -3142.5  # M67: -2517 -> -3142.5
$ws.Cells.Item(67, 14).Value = -8716.5  # N67: -7383.3335 -> -8716.5

$ws.Cells.Item(69, 8).Value = 19687.5  # H69: 18933.111 -> 19687.5
$ws.Cells.Item(69, 9).Value = 19166.666  # I69: 17599.5 -> 19166.666
$ws.Cells.Item(69, 11).Value = 57499.99800000001  # K69: 52798.5 -> 57499.99800000001
$ws.Cells.Item(69, 13).Value = -56625.99800000001  # M69: -51924.5 -> -56625.99800000001

$ws.Cells.Item(70, 8).Value = 5011.524  # H70: 5137.15 -> 5011.524
$ws.Cells.Item(70, 10).Value = 12478  # J70: 14972.75 -> 12478
$ws.Cells.Item(70, 12).Value = 37434  # L70: 44918.25 -> 37434
$ws.Cells.Item(70, 14).Value = -37974  # N70: -45458.25 -> -37974

$ws.Cells.Item(72, 8).Value = 19687.5  # H72: 18933.111 -> 19687.5
$ws.Cells.Item(72, 9).Value = 19166.666  # I72: 17599.5 -> 19166.666
$ws.Cells.Item(72, 11).Value = 172499.994  # K72: 158395.5 -> 172499.994
$ws.Cells.Item(72, 13).Value = -168131.994  # M72: -154027.5 -> -168131.994

$ws.Cells.Item(73, 8).Value = 5011.524  # H73: 5137.15 -> 5011.524
$ws.Cells.Item(73, 10).Value = 12478  # J73: 14972.75 -> 12478
$ws.Cells.Item(73, 12).Value = 37434  # L73: 44918.25 -> 37434
$ws.Cells.Item(73, 14).Value = -39306  # N73: -46790.25 -> -39306

$ws.Cells.Item(76, 8).Value = 0  # H76: 1800 -> 0
$ws.Cells.Item(76, 9).Value = 0  # I76: 1800 -> 0
$ws.Cells.Item(76, 11).Value = 0  # K76: 1800 -> 0
$ws.Cells.Item(76, 13).Value = ""  # M76: remove (was -1485)

$ws.Cells.Item(79, 8).Value = 0  # H79: 1800 -> 0
$ws.Cells.Item(79, 9).Value = 0  # I79: 1800 -> 0
$ws.Cells.Item(79, 11).Value = 0  # K79: 1800 -> 0
$ws.Cells.Item(79, 13).Value = ""  # M79: remove (was -708)

$ws.Cells.Item(86, 8).Value = 8433  # H86: 8599.666999999999 -> 8433
$ws.Cells.Item(86, 9).Value = 5500  # I86: 6000 -> 5500
$ws.Cells.Item(86, 11).Value = 5500  # K86: 6000 -> 5500
$ws.Cells.Item(86, 13).Value = -4377  # M86: -4877 -> -4377

$ws.Cells.Item(87, 8).Value = 76665.664  # H87: 78570.71000000001 -> 76665.664
$ws.Cells.Item(87, 10).Value = 76665.664  # J87: 78570.71000000001 -> 76665.664
$ws.Cells.Item(87, 12).Value = 76665.664  # L87: 78570.71000000001 -> 76665.664
$ws.Cells.Item(87, 14).Value = -79161.664  # N87: -81066.71000000001 -> -79161.664

$ws.Cells.Item(89, 8).Value = 8433  # H89: 8599.666999999999 -> 8433
$ws.Cells.Item(89, 9).Value = 5500  # I89: 6000 -> 5500
$ws.Cells.Item(89, 11).Value = 27500  # K89: 30000 -> 27500
$ws.Cells.Item(89, 13).Value = -21884  # M89: -24384 -> -21884

$ws.Cells.Item(90, 8).Value = 76665.664  # H90: 78570.71000000001 -> 76665.664
$ws.Cells.Item(90, 10).Value = 76665.664  # J90: 78570.71000000001 -> 76665.664
$ws.Cells.Item(90, 12).Value = 229996.992  # L90: 235712.13 -> 229996.992
$ws.Cells.Item(90, 14).Value = -242476.992  # N90: -248192.13 -> -242476.992

$ws.Cells.Item(98, 8).Value = 998  # H98: 998.5 -> 998
$ws.Cells.Item(98, 9).Value = 0  # I98: 999 -> 0
$ws.Cells.Item(98, 11).Value = 0  # K98: 999 -> 0
$ws.Cells.Item(98, 13).Value = ""  # M98: remove (was 499)

$ws.Cells.Item(106, 8).Value = 2952.6667  # H106: 2939.25 -> 2952.6667
$ws.Cells.Item(106, 9).Value = 2934  # I106: 2922.3333 -> 2934
$ws.Cells.Item(106, 11).Value = 2934  # K106: 2922.3333 -> 2934
$ws.Cells.Item(106, 13).Value = -2303  # M106: -2291.3333 -> -2303

$ws.Cells.Item(107, 8).Value = 1288  # H107: 1269.5714 -> 1288
$ws.Cells.Item(107, 9).Value = 1288.6154  # I107: 1266.8572 -> 1288.6154
$ws.Cells.Item(107, 10).Value = 1286.6666  # J107: 1275 -> 1286.6666
$ws.Cells.Item(107, 11).Value = 1288.6154  # K107: 1266.8572 -> 1288.6154
$ws.Cells.Item(107, 12).Value = 1286.6666  # L107: 1275 -> 1286.6666
$ws.Cells.Item(107, 13).Value = 631.3846000000001  # M107: 653.1428000000001 -> 631.3846000000001
$ws.Cells.Item(107, 14).Value = -5126.6666  # N107: -5115 -> -5126.6666

$ws.Cells.Item(122, 8).Value = 998  # H122: 998.5 -> 998
$ws.Cells.Item(122, 9).Value = 0  # I122: 999 -> 0
$ws.Cells.Item(122, 11).Value = 0  # K122: 2997 -> 0
$ws.Cells.Item(122, 13).Value = ""  # M122: remove (was -547)

$ws.Cells.Item(129, 8).Value = 2574.2  # H129: 2586.7144 -> 2574.2
$ws.Cells.Item(129, 10).Value = 2716.9  # J129: 2752.2222 -> 2716.9
$ws.Cells.Item(129, 12).Value = 8150.700000000001  # L129: 8256.6666 -> 8150.700000000001
$ws.Cells.Item(129, 14).Value = -18150.7  # N129: -18256.6666 -> -18150.7

$ws.Cells.Item(137, 8).Value = 3206.963  # H137: 3210.074 -> 3206.963
$ws.Cells.Item(137, 9).Value = 1964.909  # I137: 1972.5454 -> 1964.909
$ws.Cells.Item(137, 11).Value = 5894.727000000001  # K137: 5917.6362 -> 5894.727000000001
$ws.Cells.Item(137, 13).Value = -3344.727000000001  # M137: -3367.6362 -> -3344.727000000001

$ws.Cells.Item(138, 8).Value = 2842.3684  # H138: 3180.3333 -> 2842.3684
$ws.Cells.Item(138, 10).Value = 3093.6667  # J138: 3481.5881 -> 3093.6667
$ws.Cells.Item(138, 12).Value = 9281.000100000001  # L138: 10444.7643 -> 9281.000100000001
$ws.Cells.Item(138, 14).Value = -19561.0001  # N138: -20724.7643 -> -19561.0001

$ws = $wb.Worksheets.Item(2)  # ARM
$ws.Cells.Item(61, 8).Value = 2678.6924  # H61: 2672.5715 -> 2678.6924
$ws.Cells.Item(61, 9).Value = 2485.5833  # I61: 2620.0908 -> 2485.5833
$ws.Cells.Item(61, 10).Value = 4996  # J61: 2865 -> 4996
$ws.Cells.Item(61, 11).Value = 2485.5833  # K61: 2620.0908 -> 2485.5833
$ws.Cells.Item(61, 12).Value = 4996  # L61: 2865 -> 4996
$ws.Cells.Item(61, 13).Value = -2273.5833  # M61: -2408.0908 -> -2273.5833
$ws.Cells.Item(61, 14).Value = -5420  # N61: -3289 -> -5420

$ws.Cells.Item(136, 8).Value = 2678.6924  # H136: 2672.5715 -> 2678.6924
$ws.Cells.Item(136, 9).Value = 2485.5833  # I136: 2620.0908 -> 2485.5833
$ws.Cells.Item(136, 10).Value = 4996  # J136: 2865 -> 4996
$ws.Cells.Item(136, 11).Value = 7456.749899999999  # K136: 7860.2724 -> 7456.749899999999
$ws.Cells.Item(136, 12).Value = 14988  # L136: 8595 -> 14988
$ws.Cells.Item(136, 13).Value = -4906.749899999999  # M136: -5310.2724 -> -4906.749899999999
$ws.Cells.Item(136, 14).Value = -20088  # N136: -13695 -> -20088

$ws = $wb.Worksheets.Item(3)  # BSM
$ws.Cells.Item(134, 8).Value = 1335.1333  # H134: 1360.4546 -> 1335.1333
$ws.Cells.Item(134, 9).Value = 1356.7675  # I134: 1383.8096 -> 1356.7675
$ws.Cells.Item(134, 11).Value = 4070.3025  # K134: 4151.4288 -> 4070.3025
$ws.Cells.Item(134, 13).Value = -1535.3025  # M134: -1616.4288 -> -1535.3025

$ws = $wb.Worksheets.Item(4)  # CRP
$ws.Cells.Item(58, 8).Value = 2757  # H58: 2951.7144 -> 2757
$ws.Cells.Item(58, 9).Value = 2509.3333  # I58: 2732.4 -> 2509.3333
$ws.Cells.Item(58, 11).Value = 2509.3333  # K58: 2732.4 -> 2509.3333
$ws.Cells.Item(58, 13).Value = -2306.3333  # M58: -2529.4 -> -2306.3333

$ws.Cells.Item(136, 8).Value = 2757  # H136: 2951.7144 -> 2757
$ws.Cells.Item(136, 9).Value = 2509.3333  # I136: 2732.4 -> 2509.3333
$ws.Cells.Item(136, 11).Value = 7527.999899999999  # K136: 8197.200000000001 -> 7527.999899999999
$ws.Cells.Item(136, 13).Value = -4977.999899999999  # M136: -5647.200000000001 -> -4977.999899999999

$ws = $wb.Worksheets.Item(5)  # CUL
$ws.Cells.Item(68, 8).Value = 4140.2  # H68: 4512.5 -> 4140.2
$ws.Cells.Item(68, 10).Value = 4140.2  # J68: 4512.5 -> 4140.2
$ws.Cells.Item(68, 12).Value = 12420.6  # L68: 13537.5 -> 12420.6
$ws.Cells.Item(68, 14).Value = -14042.6  # N68: -15159.5 -> -14042.6

$ws.Cells.Item(71, 8).Value = 4140.2  # H71: 4512.5 -> 4140.2
$ws.Cells.Item(71, 10).Value = 4140.2  # J71: 4512.5 -> 4140.2
$ws.Cells.Item(71, 12).Value = 37261.8  # L71: 40612.5 -> 37261.8
$ws.Cells.Item(71, 14).Value = -45373.8  # N71: -48724.5 -> -45373.8

$ws.Cells.Item(98, 8).Value = 405.875  # H98: 587.125 -> 405.875
$ws.Cells.Item(98, 9).Value = 249.57143  # I98: 283 -> 249.57143
$ws.Cells.Item(98, 10).Value = 1500  # J98: 1499.5 -> 1500
$ws.Cells.Item(98, 11).Value = 748.71429  # K98: 849 -> 748.71429
$ws.Cells.Item(98, 12).Value = 4500  # L98: 4498.5 -> 4500
$ws.Cells.Item(98, 13).Value = 749.28571  # M98: 649 -> 749.28571
$ws.Cells.Item(98, 14).Value = -7496  # N98: -7494.5 -> -7496

$ws.Cells.Item(117, 8).Value = 2832.3333  # H117: 2078 -> 2832.3333
$ws.Cells.Item(117, 9).Value = 2832.3333  # I117: 2267 -> 2832.3333
$ws.Cells.Item(117, 10).Value = 0  # J117: 1700 -> 0
$ws.Cells.Item(117, 11).Value = 8496.999899999999  # K117: 6801 -> 8496.999899999999
$ws.Cells.Item(117, 12).Value = 0  # L117: 5100 -> 0
$ws.Cells.Item(117, 13).Value = -5054.999899999999  # M117: -3359 -> -5054.999899999999
$ws.Cells.Item(117, 14).Value = ""  # N117: remove (was -11984)

$ws.Cells.Item(120, 8).Value = 14999  # H120: 14999.5 -> 14999
$ws.Cells.Item(120, 9).Value = 14998  # I120: 0 -> 14998
$ws.Cells.Item(120, 10).Value = 15000  # J120: 14999.5 -> 15000
$ws.Cells.Item(120, 11).Value = 44994  # K120: 0 -> 44994
$ws.Cells.Item(120, 12).Value = 45000  # L120: 44998.5 -> 45000
$ws.Cells.Item(120, 13).Value = -40156  # M120: None -> -40156
$ws.Cells.Item(120, 14).Value = -54676  # N120: -54674.5 -> -54676

$ws.Cells.Item(122, 8).Value = 71430030  # H122: 166668160 -> 71430030
$ws.Cells.Item(122, 9).Value = 1574.75  # I122: 2000 -> 1574.75
$ws.Cells.Item(122, 10).Value = 166667970  # J122: 250001250 -> 166667970
$ws.Cells.Item(122, 11).Value = 14172.75  # K122: 18000 -> 14172.75
$ws.Cells.Item(122, 12).Value = 1500011730  # L122: 2250011250 -> 1500011730
$ws.Cells.Item(122, 13).Value = -11722.75  # M122: -15550 -> -11722.75
$ws.Cells.Item(122, 14).Value = -1500016630  # N122: -2250016150 -> -1500016630

$ws = $wb.Worksheets.Item(6)  # GSM
$ws.Cells.Item(2, 8).Value = 1021.3333  # H2: 1195.6 -> 1021.3333
$ws.Cells.Item(2, 10).Value = 2075  # J2: 4000 -> 2075
$ws.Cells.Item(2, 12).Value = 2075  # L2: 4000 -> 2075
$ws.Cells.Item(2, 14).Value = -2301  # N2: -4226 -> -2301

$ws.Cells.Item(80, 8).Value = 0  # H80: 800 -> 0
$ws.Cells.Item(80, 9).Value = 0  # I80: 800 -> 0
$ws.Cells.Item(80, 11).Value = 0  # K80: 800 -> 0
$ws.Cells.Item(80, 13).Value = ""  # M80: remove (was 198)

$ws.Cells.Item(83, 8).Value = 0  # H83: 800 -> 0
$ws.Cells.Item(83, 9).Value = 0  # I83: 800 -> 0
$ws.Cells.Item(83, 11).Value = 0  # K83: 4000 -> 0
$ws.Cells.Item(83, 13).Value = ""  # M83: remove (was 992)

$ws.Cells.Item(128, 8).Value = 59989.5  # H128: 59992 -> 59989.5
$ws.Cells.Item(128, 10).Value = 59989.5  # J128: 59992 -> 59989.5
$ws.Cells.Item(128, 12).Value = 59989.5  # L128: 59992 -> 59989.5
$ws.Cells.Item(128, 14).Value = -69949.5  # N128: -69952 -> -69949.5

$ws = $wb.Worksheets.Item(7)  # LTW
$ws.Cells.Item(68, 8).Value = 2248  # H68: 2160.3333 -> 2248
$ws.Cells.Item(68, 9).Value = 2248  # I68: 2160.3333 -> 2248
$ws.Cells.Item(68, 11).Value = 2248  # K68: 2160.3333 -> 2248
$ws.Cells.Item(68, 13).Value = -1499  # M68: -1411.3333 -> -1499

$ws.Cells.Item(71, 8).Value = 2248  # H71: 2160.3333 -> 2248
$ws.Cells.Item(71, 9).Value = 2248  # I71: 2160.3333 -> 2248
$ws.Cells.Item(71, 11).Value = 11240  # K71: 10801.6665 -> 11240
$ws.Cells.Item(71, 13).Value = -7496  # M71: -7057.666499999999 -> -7496

$ws.Cells.Item(82, 8).Value = 3822.182  # H82: 3246.7693 -> 3822.182
$ws.Cells.Item(82, 9).Value = 2210  # I82: 1905.7142 -> 2210
$ws.Cells.Item(82, 10).Value = 5756.8  # J82: 4811.3335 -> 5756.8
$ws.Cells.Item(82, 11).Value = 2210  # K82: 1905.7142 -> 2210
$ws.Cells.Item(82, 12).Value = 5756.8  # L82: 4811.3335 -> 5756.8
$ws.Cells.Item(82, 13).Value = -1849  # M82: -1544.7142 -> -1849
$ws.Cells.Item(82, 14).Value = -6478.8  # N82: -5533.3335 -> -6478.8

$ws.Cells.Item(85, 8).Value = 3822.182  # H85: 3246.7693 -> 3822.182
$ws.Cells.Item(85, 9).Value = 2210  # I85: 1905.7142 -> 2210
$ws.Cells.Item(85, 10).Value = 5756.8  # J85: 4811.3335 -> 5756.8
$ws.Cells.Item(85, 11).Value = 2210  # K85: 1905.7142 -> 2210
$ws.Cells.Item(85, 12).Value = 5756.8  # L85: 4811.3335 -> 5756.8
$ws.Cells.Item(85, 13).Value = -962  # M85: -657.7141999999999 -> -962
$ws.Cells.Item(85, 14).Value = -8252.799999999999  # N85: -7307.3335 -> -8252.799999999999

$ws = $wb.Worksheets.Item(8)  # WVR
$ws.Cells.Item(62, 8).Value = 6333.3335  # H62: 7750 -> 6333.3335
$ws.Cells.Item(62, 9).Value = 5000  # I62: 7500 -> 5000
$ws.Cells.Item(62, 10).Value = 7000  # J62: 8000 -> 7000
$ws.Cells.Item(62, 11).Value = 5000  # K62: 7500 -> 5000
$ws.Cells.Item(62, 12).Value = 7000  # L62: 8000 -> 7000
$ws.Cells.Item(62, 13).Value = -4376  # M62: -6876 -> -4376
$ws.Cells.Item(62, 14).Value = -8248  # N62: -9248 -> -8248

$ws.Cells.Item(65, 8).Value = 6333.3335  # H65: 7750 -> 6333.3335
$ws.Cells.Item(65, 9).Value = 5000  # I65: 7500 -> 5000
$ws.Cells.Item(65, 10).Value = 7000  # J65: 8000 -> 7000
$ws.Cells.Item(65, 11).Value = 25000  # K65: 37500 -> 25000
$ws.Cells.Item(65, 12).Value = 35000  # L65: 40000 -> 35000
$ws.Cells.Item(65, 13).Value = -21880  # M65: -34380 -> -21880
$ws.Cells.Item(65, 14).Value = -41240  # N65: -46240 -> -41240

$ws.Cells.Item(70, 8).Value = 55000  # H70: 25000 -> 55000
$ws.Cells.Item(70, 9).Value = 0  # I70: 25000 -> 0
$ws.Cells.Item(70, 10).Value = 55000  # J70: 0 -> 55000
$ws.Cells.Item(70, 11).Value = 0  # K70: 25000 -> 0
$ws.Cells.Item(70, 12).Value = 55000  # L70: 0 -> 55000
$ws.Cells.Item(70, 13).Value = ""  # M70: remove (was -24685)
$ws.Cells.Item(70, 14).Value = -55630  # N70: None -> -55630

$ws.Cells.Item(73, 8).Value = 55000  # H73: 25000 -> 55000
$ws.Cells.Item(73, 9).Value = 0  # I73: 25000 -> 0
$ws.Cells.Item(73, 10).Value = 55000  # J73: 0 -> 55000
$ws.Cells.Item(73, 11).Value = 0  # K73: 25000 -> 0
$ws.Cells.Item(73, 12).Value = 55000  # L73: 0 -> 55000
$ws.Cells.Item(73, 13).Value = ""  # M73: remove (was -23908)
$ws.Cells.Item(73, 14).Value = -57184  # N73: None -> -57184
